$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.280.85"
$ws.Range("E2").Value = "  +2.64%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.899.15"
$ws.Range("E3").Value = "  +0.87%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  -1.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.50"
$ws.Range("E5").Value = "  -0.42%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.005"
$ws.Range("E6").Value = "  -0.93%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5142"
$ws.Range("E7").Value = "  +0.26%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3923"
$ws.Range("E8").Value = "  -1.38%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08440"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.53"
$ws.Range("E10").Value = "  +1.56%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.116"
$ws.Range("E11").Value = "  +0.34%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.246"
$ws.Range("E12").Value = "  -0.62%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.895.83"
$ws.Range("E13").Value = "  +0.78%  "
$ws.Range("E14").Value = "  +0.68%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.326"
$ws.Range("E15").Value = "  +0.56%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.005"
$ws.Range("E16").Value = "  -1.20%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "93.34"
$ws.Range("E17").Value = "  +2.03%  "
$ws.Range("E18").Value = "  -0.32%  "
$ws.Range("E19").Value = "  -0.24%  "
$ws.Range("E20").Value = "  +0.42%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.004"
$ws.Range("E21").Value = "  -0.80%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.024"
$ws.Range("E22").Value = "  +1.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "29.282.37"
$ws.Range("E23").Value = "  +2.36%  "
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.218"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.114.34"
$ws.Range("E26").Value = "  +0.83%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.97"
$ws.Range("E27").Value = "  +0.53%  "
$ws.Range("B28").Value = "Monero"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "159.07"
$ws.Range("E28").Value = "  -1.54%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.445"
$ws.Range("E29").Value = "  +2.36%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "128.45"
$ws.Range("E30").Value = "  +0.83%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.059"
$ws.Range("E31").Value = "  +0.82%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1046"
$ws.Range("E32").Value = "  -1.09%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.160"
$ws.Range("E33").Value = "  +6.27%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.663"
$ws.Range("E34").Value = "  +1.12%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02475"
$ws.Range("E35").Value = "  +1.47%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06570"
$ws.Range("E36").Value = "  +0.95%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "9.055"
$ws.Range("E37").Value = "  +1.52%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2194"
$ws.Range("E38").Value = "  +0.28%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.234"
$ws.Range("E39").Value = "  +3.47%  "
$ws.Range("E40").Value = "  +1.47%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6500"
$ws.Range("E41").Value = "  +0.55%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.233"
$ws.Range("E42").Value = "  -2.58%  "
$ws.Range("E43").Value = "  +0.33%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6054"
$ws.Range("E44").Value = "  -0.59%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.22"
$ws.Range("E45").Value = "  +1.40%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.675"
$ws.Range("E46").Value = "  -1.10%  "
$ws.Range("E47").Value = "  +1.76%  "
$ws.Range("E48").Value = "  +1.31%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "123.39"
$ws.Range("E49").Value = "  +0.53%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.177"
$ws.Range("E50").Value = "  -2.50%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "77.69"
$ws.Range("E51").Value = "  +0.72%  "
